# Apply the "25hike" scenario update: shock input changes from 1.5 to 1.25
# and all downstream computed columns (A-E) are refreshed with new results.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 2).Value = 1.25
$ws.Cells.Item(1, 3).Value = 1.25
$ws.Cells.Item(1, 4).Value = 1.25

$ws.Cells.Item(2, 1).Value = 0.049539478305621484
$ws.Cells.Item(2, 2).Value = 1.482828648036421
$ws.Cells.Item(2, 3).Value = 1.4273375684690084
$ws.Cells.Item(2, 4).Value = 1.2709545783232779
$ws.Cells.Item(2, 5).Value = -0.00001000000000000119

$ws.Cells.Item(3, 1).Value = 0.13414609359043725
$ws.Cells.Item(3, 2).Value = 1.675780476516227
$ws.Cells.Item(3, 3).Value = 1.5472268698337288
$ws.Cells.Item(3, 4).Value = 1.3073889091606432
$ws.Cells.Item(3, 5).Value = -0.000010000000000000739

$ws.Cells.Item(4, 1).Value = 0.24421067667960472
$ws.Cells.Item(4, 2).Value = 1.8579682183076223
$ws.Cells.Item(4, 3).Value = 1.6535967872670894
$ws.Cells.Item(4, 4).Value = 1.3569410469838714
$ws.Cells.Item(4, 5).Value = -0.000009999999999995127

$ws.Cells.Item(5, 1).Value = 0.32953850322660372
$ws.Cells.Item(5, 2).Value = 2.0059462753145656
$ws.Cells.Item(5, 3).Value = 1.7416207125319034
$ws.Cells.Item(5, 4).Value = 1.415351517533634
$ws.Cells.Item(5, 5).Value = 0.25

$ws.Cells.Item(6, 1).Value = 0.36708716631807436
$ws.Cells.Item(6, 2).Value = 2.1000997260942427
$ws.Cells.Item(6, 3).Value = 1.8049776316410207
$ws.Cells.Item(6, 4).Value = 1.4769788563040886
$ws.Cells.Item(6, 5).Value = 0.50441327106219147

$ws.Cells.Item(7, 1).Value = 0.35650182343746345
$ws.Cells.Item(7, 2).Value = 2.1381291896590771
$ws.Cells.Item(7, 3).Value = 1.8444747356646456
$ws.Cells.Item(7, 4).Value = 1.5364823863060377
$ws.Cells.Item(7, 5).Value = 0.69892739665131642

$ws.Cells.Item(8, 1).Value = 0.31620558781400043
$ws.Cells.Item(8, 2).Value = 2.1349086255526699
$ws.Cells.Item(8, 3).Value = 1.8658800774563511
$ws.Cells.Item(8, 4).Value = 1.5903407478382345
$ws.Cells.Item(8, 5).Value = 0.78030358720716331

$ws.Cells.Item(9, 1).Value = 0.2661805135780213
$ws.Cells.Item(9, 2).Value = 2.1106736723741228
$ws.Cells.Item(9, 3).Value = 1.8768310834529216
$ws.Cells.Item(9, 4).Value = 1.6371707110464644
$ws.Cells.Item(9, 5).Value = 0.78511296564742927

$ws.Cells.Item(10, 1).Value = 0.22016011688356751
$ws.Cells.Item(10, 2).Value = 2.0820698089607355
$ws.Cells.Item(10, 3).Value = 1.8837549432595702
$ws.Cells.Item(10, 4).Value = 1.6772116298587489
$ws.Cells.Item(10, 5).Value = 0.7546908417544429

$ws.Cells.Item(11, 1).Value = 0.18422664229325839
$ws.Cells.Item(11, 2).Value = 2.0586804061580057
$ws.Cells.Item(11, 3).Value = 1.8905381232897964
$ws.Cells.Item(11, 4).Value = 1.711543819725682
$ws.Cells.Item(11, 5).Value = 0.71970037087417282

$ws.Cells.Item(12, 1).Value = 0.15872259311044459
$ws.Cells.Item(12, 2).Value = 2.0435928220934416
$ws.Cells.Item(12, 3).Value = 1.898654848834715
$ws.Cells.Item(12, 4).Value = 1.7414282299387804
$ws.Cells.Item(12, 5).Value = 0.69593118327704817

$ws.Cells.Item(13, 1).Value = 0.14106041591928237
$ws.Cells.Item(13, 2).Value = 2.0357708092675408
$ws.Cells.Item(13, 3).Value = 1.9079960058588568
$ws.Cells.Item(13, 4).Value = 1.7679190620783689
$ws.Cells.Item(13, 5).Value = 0.68702652005285147

$ws.Cells.Item(14, 1).Value = 0.12802841753733576
$ws.Cells.Item(14, 2).Value = 2.0325280789108744
$ws.Cells.Item(14, 3).Value = 1.9177727558025737
$ws.Cells.Item(14, 4).Value = 1.7917338735932944
$ws.Cells.Item(14, 5).Value = 0.68969543000576261

$ws.Cells.Item(15, 1).Value = 0.11710118413182996
$ws.Cells.Item(15, 2).Value = 2.0312466290551487
$ws.Cells.Item(15, 3).Value = 1.9271640785754036
$ws.Cells.Item(15, 4).Value = 1.8132900215848613
$ws.Cells.Item(15, 5).Value = 0.69843759197845157

$ws.Cells.Item(16, 1).Value = 0.10684755665791355
$ws.Cells.Item(16, 2).Value = 2.0301758284179581
$ws.Cells.Item(16, 3).Value = 1.9356322630129754
$ws.Cells.Item(16, 4).Value = 1.8328097441998401
$ws.Cells.Item(16, 5).Value = 0.70850734603149168

$ws.Cells.Item(17, 1).Value = 0.096774989132813696
$ws.Cells.Item(17, 2).Value = 2.0285331195622223
$ws.Cells.Item(17, 3).Value = 1.9429778864962757
$ws.Cells.Item(17, 4).Value = 1.8504248479824543
$ws.Cells.Item(17, 5).Value = 0.71707764509176086

$ws.Cells.Item(18, 1).Value = 0.086955304163263303
$ws.Cells.Item(18, 2).Value = 2.0262388670418998
$ws.Cells.Item(18, 3).Value = 1.9492530707801157
$ws.Cells.Item(18, 4).Value = 1.8662481096978045
$ws.Cells.Item(18, 5).Value = 0.72316885295556466

$ws.Cells.Item(19, 1).Value = 0.077667424784810146
$ws.Cells.Item(19, 2).Value = 2.0235628601329849
$ws.Cells.Item(19, 3).Value = 1.9546356951393964
$ws.Cells.Item(19, 4).Value = 1.8804064372369707
$ws.Cells.Item(19, 5).Value = 0.72701260424075986

$ws.Cells.Item(20, 1).Value = 0.069163872027702356
$ws.Cells.Item(20, 2).Value = 2.0208434802619761
$ws.Cells.Item(20, 3).Value = 1.9593261644500033
$ws.Cells.Item(20, 4).Value = 1.8930457711092212
$ws.Cells.Item(20, 5).Value = 0.72935354739595382

$ws.Cells.Item(21, 1).Value = 0.0615717186971333
$ws.Cells.Item(21, 2).Value = 2.0183326407871034
$ws.Cells.Item(21, 3).Value = 1.9634887322974461
$ws.Cells.Item(21, 4).Value = 1.9043215893802308
$ws.Cells.Item(21, 5).Value = 0.73094865524743424

$ws.Cells.Item(22, 1).Value = 0.054889424530941082
$ws.Cells.Item(22, 2).Value = 2.0161509447354198
$ws.Cells.Item(22, 3).Value = 1.967233332515212
$ws.Cells.Item(22, 4).Value = 1.9143862313621978
$ws.Cells.Item(22, 5).Value = 0.73232356973955259

$ws.Cells.Item(23, 1).Value = 0.049030038327224051
$ws.Cells.Item(23, 2).Value = 2.0143097896778306
$ws.Cells.Item(23, 3).Value = 1.9706228237015808
$ws.Cells.Item(23, 4).Value = 1.9233793516106048
$ws.Cells.Item(23, 5).Value = 0.73373069930579304

$ws.Cells.Item(24, 1).Value = 0.043872143635534287
$ws.Cells.Item(24, 2).Value = 2.0127581049117769
$ws.Cells.Item(24, 3).Value = 1.9736901091979784
$ws.Cells.Item(24, 4).Value = 1.9314234394077103
$ws.Cells.Item(24, 5).Value = 0.73521715443973468

$ws.Cells.Item(25, 1).Value = 0.03929765065022079
$ws.Cells.Item(25, 2).Value = 2.0114254305476753
$ws.Cells.Item(25, 3).Value = 1.9764545937329712
$ws.Cells.Item(25, 4).Value = 1.9386236186103072
$ws.Cells.Item(25, 5).Value = 0.73672291949507118

$ws.Cells.Item(26, 1).Value = 0.035210803425341317
$ws.Cells.Item(26, 2).Value = 2.0102492605353239
$ws.Cells.Item(26, 3).Value = 1.9789329473346391
$ws.Cells.Item(26, 4).Value = 1.9450699263835585
$ws.Cells.Item(26, 5).Value = 0.73816161461319096

$ws.Cells.Item(27, 1).Value = 0.031541815945485621
$ws.Cells.Item(27, 2).Value = 2.0091855800358212
$ws.Cells.Item(27, 3).Value = 1.9811439523874685
$ws.Cells.Item(27, 4).Value = 1.9508403352122623
$ws.Cells.Item(27, 5).Value = 0.73946650821361415

$ws.Cells.Item(28, 1).Value = 0.028241852517580029
$ws.Cells.Item(28, 2).Value = 2.0082082838412902
$ws.Cells.Item(28, 3).Value = 1.9831085396237667
$ws.Cells.Item(28, 4).Value = 1.9560034505888746
$ws.Cells.Item(28, 5).Value = 0.74060524720019061

$ws.Cells.Item(29, 1).Value = 0.025275259183865414
$ws.Cells.Item(29, 2).Value = 2.00730232245829
$ws.Cells.Item(29, 3).Value = 1.9848482849570153
$ws.Cells.Item(29, 4).Value = 1.9606203490571221
$ws.Cells.Item(29, 5).Value = 0.741574152432811

$ws.Cells.Item(30, 1).Value = 0.022613010035249691
$ws.Cells.Item(30, 2).Value = 2.0064574947191951
$ws.Cells.Item(30, 3).Value = 1.9863820424993124
$ws.Cells.Item(30, 4).Value = 1.9647456921667086
$ws.Cells.Item(30, 5).Value = 0.74238546201364985

$ws.Cells.Item(31, 1).Value = 0.020228437025464344
$ws.Cells.Item(31, 2).Value = 2.0056609191219845
$ws.Cells.Item(31, 3).Value = 1.9877253465759206
$ws.Cells.Item(31, 4).Value = 1.9684280625926838
$ws.Cells.Item(31, 5).Value = 0.74305292330189154

$ws.Cells.Item(32, 1).Value = 0.0180958470535597
$ws.Cells.Item(32, 2).Value = 2.0048966742227448
$ws.Cells.Item(32, 3).Value = 1.9888857551083681
$ws.Cells.Item(32, 4).Value = 1.9717102376393891
$ws.Cells.Item(32, 5).Value = 0.74358442795806412

$ws.Cells.Item(33, 1).Value = 0.0161898600680737
$ws.Cells.Item(33, 2).Value = 2.0041369798466797
$ws.Cells.Item(33, 3).Value = 1.9898668332535303
$ws.Cells.Item(33, 4).Value = 1.9746286444380452
$ws.Cells.Item(33, 5).Value = 0.74397243561096715

$ws.Cells.Item(34, 1).Value = 0.014487059723701106
$ws.Cells.Item(34, 2).Value = 2.003351606490861
$ws.Cells.Item(34, 3).Value = 1.9906560713036345
$ws.Cells.Item(34, 4).Value = 1.9772137110227985
$ws.Cells.Item(34, 5).Value = 0.74419812083014347

$ws.Cells.Item(35, 1).Value = 0.01296468156196476
$ws.Cells.Item(35, 2).Value = 2.0024826097544777
$ws.Cells.Item(35, 3).Value = 1.9912411439030691
$ws.Cells.Item(35, 4).Value = 1.9794879119086497
$ws.Cells.Item(35, 5).Value = 0.74421434396775121

$ws.Cells.Item(36, 1).Value = 0.011604649924566137
$ws.Cells.Item(36, 2).Value = 2.0014796139281077
$ws.Cells.Item(36, 3).Value = 1.9915701305327689
$ws.Cells.Item(36, 4).Value = 1.981467165090401
$ws.Cells.Item(36, 5).Value = 0.7439662509438455

$ws.Cells.Item(37, 1).Value = 0.010386870474370424
$ws.Cells.Item(37, 2).Value = 2.000215436391831
$ws.Cells.Item(37, 3).Value = 1.9916107260093663
$ws.Cells.Item(37, 4).Value = 1.9831545095075298
$ws.Cells.Item(37, 5).Value = 0.74333451628941116

$ws.Cells.Item(38, 1).Value = 0.00930171426906958
$ws.Cells.Item(38, 2).Value = 1.9986140884837003
$ws.Cells.Item(38, 3).Value = 1.9912175487079276
$ws.Cells.Item(38, 4).Value = 1.9845458716153852
$ws.Cells.Item(38, 5).Value = 0.74221018834425334

$ws.Cells.Item(39, 1).Value = 0.0083262888796092529
$ws.Cells.Item(39, 2).Value = 1.9963691557686407
$ws.Cells.Item(39, 3).Value = 1.9903454249442811
$ws.Cells.Item(39, 4).Value = 1.9856099671891783
$ws.Cells.Item(39, 5).Value = 0.74030708264359502

$ws.Cells.Item(40, 1).Value = 0.0074679310721769993
$ws.Cells.Item(40, 2).Value = 1.9934037846973747
$ws.Cells.Item(40, 3).Value = 1.9886029950233515
$ws.Cells.Item(40, 4).Value = 1.986311410764916
$ws.Cells.Item(40, 5).Value = 0.73743803192739121

$ws.Cells.Item(41, 1).Value = 0.0066827401243234105
$ws.Cells.Item(41, 2).Value = 1.9889212506159581
$ws.Cells.Item(41, 3).Value = 1.9860105869690967
$ws.Cells.Item(41, 4).Value = 1.9865462963515097
$ws.Cells.Item(41, 5).Value = 0.73289354764122905

$ws.Cells.Item(42, 1).Value = 0.006026886150418455
$ws.Cells.Item(42, 2).Value = 1.983034868422328
$ws.Cells.Item(42, 3).Value = 1.9814836853480449
$ws.Cells.Item(42, 4).Value = 1.9862302678378834
$ws.Cells.Item(42, 5).Value = 0.72644087132827395

$ws.Cells.Item(43, 1).Value = 0.0053751687410549904
$ws.Cells.Item(43, 2).Value = 1.9735358422263429
$ws.Cells.Item(43, 3).Value = 1.9754988760683285
$ws.Cells.Item(43, 4).Value = 1.9850877695328446
$ws.Cells.Item(43, 5).Value = 0.71623582588922607

$ws.Cells.Item(44, 1).Value = 0.0049451578629795901
$ws.Cells.Item(44, 2).Value = 1.9616154254435942
$ws.Cells.Item(44, 3).Value = 1.9649057080006238
$ws.Cells.Item(44, 4).Value = 1.9829752585648122
$ws.Cells.Item(44, 5).Value = 0.70237303362223469

$ws.Cells.Item(45, 1).Value = 0.0043229360377254967
$ws.Cells.Item(45, 2).Value = 1.9407917972969608
$ws.Cells.Item(45, 3).Value = 1.9522464144205207
$ws.Cells.Item(45, 4).Value = 1.9791787470507054
$ws.Cells.Item(45, 5).Value = 0.67979056212304667

$ws.Cells.Item(46, 1).Value = 0.0042773876684296134
$ws.Cells.Item(46, 2).Value = 1.9168891380939637
$ws.Cells.Item(46, 3).Value = 1.9279277843888702
$ws.Cells.Item(46, 4).Value = 1.9735726822445989
$ws.Cells.Item(46, 5).Value = 0.65075525515231125

$ws.Cells.Item(47, 1).Value = 0.0033869363573466836
$ws.Cells.Item(47, 2).Value = 1.8699520478939862
$ws.Cells.Item(47, 3).Value = 1.9026787591849059
$ws.Cells.Item(47, 4).Value = 1.9642468251530438
$ws.Cells.Item(47, 5).Value = 0.60058447957917582

$ws.Cells.Item(48, 1).Value = 0.004319512419410477
$ws.Cells.Item(48, 2).Value = 1.8237734230919034
$ws.Cells.Item(48, 3).Value = 1.8461155517749903
$ws.Cells.Item(48, 4).Value = 1.9516042189675411
$ws.Cells.Item(48, 5).Value = 0.54128547955764605

$ws.Cells.Item(49, 1).Value = 0.0021534025700575103
$ws.Cells.Item(49, 2).Value = 1.71457802133745
$ws.Cells.Item(49, 3).Value = 1.7995208134636203
$ws.Cells.Item(49, 4).Value = 1.9302718611808327
$ws.Cells.Item(49, 5).Value = 0.42811694889043039

$ws.Cells.Item(50, 1).Value = 0.006115616941848318
$ws.Cells.Item(50, 2).Value = 1.6324954887812864
$ws.Cells.Item(50, 3).Value = 1.663771528608301
$ws.Cells.Item(50, 4).Value = 1.9034719876648736
$ws.Cells.Item(50, 5).Value = 0.31155480082003573

$ws.Cells.Item(51, 1).Value = -0.00085290451124724899
$ws.Cells.Item(51, 2).Value = 1.3680959749959003
$ws.Cells.Item(51, 3).Value = 1.5906544299699157
$ws.Cells.Item(51, 4).Value = 1.855288146524666
$ws.Cells.Item(51, 5).Value = 0.049744166436723747

$ws.Cells.Item(52, 2).Value = 1.25
$ws.Cells.Item(52, 3).Value = 1.25
$ws.Cells.Item(52, 4).Value = 1.25
